$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values. "D" column (Price) holds numeric-looking text,
# so we force text format on those cells to preserve the original string type.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '239.32'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '21.43'

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.130'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05552'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.366'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.377'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8068'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9404'

# Row 10
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1397'
$ws.Range('E10').Value = '9WazirXWRX'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07317'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03087'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03058'
$ws.Range('E13').Value = '12BitrueCoinBTR'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09292'
$ws.Range('E14').Value = '13BitMartTokenBMX'

# Row 15
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.630'
$ws.Range('E15').Value = '14MCDexMCB'

# Row 16
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001648'
$ws.Range('E16').Value = '15BitForexTokenBF'

# Row 17
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04711'
$ws.Range('E17').Value = '16CoinExTokenCET'

# Row 18
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0005756'
$ws.Range('E18').Value = '17OneONE'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006377'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.004946'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.001054'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001511'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0003108'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.761'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.101'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.3251'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1286'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03870'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006897'

# Row 42
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1035'
$ws.Range('E42').Value = '41BKEXTokenBKK'

# Row 43
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003107'
$ws.Range('E43').Value = '42CEJICEJI'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008357'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005949'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000752'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0005516'

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.6842'

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1019'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002104'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.01012'
